# Update the ifo QoQ GVA matched error table with the latest preprocessing
# results (added ifo GDP component analysis preprocessing).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = @(0.07460170201535693, 0.4066704939480233, 0.2229059984459623, 0.472129217954113, 0.4825607857740096, 15)
    3  = @(0.314504379002312, 0.4790566883332484, 0.3549208902847111, 0.5957523732933937, 0.5250721856764178, 14)
    4  = @(0.4632337633744819, 0.5692069213481521, 0.4862027082318384, 0.6972823733838669, 0.5424499556889087, 13)
    5  = @(0.6092340853067807, 0.6750741196220317, 0.6306383360763419, 0.7941274054434476, 0.5320343448005646, 12)
    6  = @(0.638688239432181, 0.7242386063103117, 0.7271270905885431, 0.8527174740724756, 0.5925579007485432, 11)
    7  = @(0.5167424503166, 0.6241159593079864, 0.538927811949409, 0.7341170287831559, 0.5496514572286344, 10)
    8  = @(0.370868480679851, 0.5185686719842275, 0.3627064250084273, 0.6022511311806955, 0.5032974959479504, 9)
    9  = @(0.2966541511547854, 0.4816046978711712, 0.2612550586842152, 0.5111311560492231, 0.4559623317163587, 6)
    10 = @(0.4961513153058069, 0.4961513153058069, 0.3214734301341839, 0.5669862697933556, 0.3360966433657923, 3)
}

foreach ($row in $values.Keys) {
    $rowVals = $values[$row]
    $ws.Cells.Item($row, 2).Value = $rowVals[0]
    $ws.Cells.Item($row, 3).Value = $rowVals[1]
    $ws.Cells.Item($row, 4).Value = $rowVals[2]
    $ws.Cells.Item($row, 5).Value = $rowVals[3]
    $ws.Cells.Item($row, 6).Value = $rowVals[4]
    $ws.Cells.Item($row, 7).Value = $rowVals[5]
}
